# Se agregó un quinto sensor de temperatura (DHT22) y se agregaron datos
# faltantes de los sensores anteriores a la tabla.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# ---------------------------------------------------------------------
# 1) Datos faltantes del sensor TMP36 (fila 5)
# ---------------------------------------------------------------------
$ws.Range("I5").Value = "NO"
$ws.Range("L5").Value = 3
$ws.Range("M5").Value = "Vin, Vout, GND"

# ---------------------------------------------------------------------
# 2) Datos faltantes del sensor TC74 (fila 6)
# ---------------------------------------------------------------------
$ws.Range("I6").Value = "NO"
$ws.Range("K6").Value = "5 µA"
$ws.Range("K6").WrapText = $false
$ws.Range("M6").Value = "NC, SDA, GND, VDD, SCLIK"

# ---------------------------------------------------------------------
# 3) Datos faltantes del sensor DHT11 (fila 7)
# ---------------------------------------------------------------------
$ws.Range("I7").Value = "SI"
$ws.Range("K7").Value = "2.5 µA"
$ws.Range("K7").WrapText = $false

# ---------------------------------------------------------------------
# 4) Nuevo quinto sensor: DHT22 (fila 8)
# ---------------------------------------------------------------------
# Hereda el formato (bordes, fuente, alineación) de la fila anterior
$ws.Range("B7:R7").Copy()
$ws.Range("B8:R8").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Rows("8").RowHeight = 75

$ws.Range("B8").Value = "Sensor temperatura  DHT22"
$ws.Range("C8").Value = "necesita medir los niveles de temperatura y humedad de ciertos ambientes"
$ws.Range("D8").Value = 20
$ws.Range("E8").Value = 15
$ws.Range("F8").Value = 8
$ws.Range("G8").Value = "2 segundos"
$ws.Range("H8").Value = "3.3 V a 6V"
$ws.Range("I8").Value = "SI"
$ws.Range("J8").Value = "0.5°C"
$ws.Range("K8").Value = "2.5 µA"
$ws.Range("L8").Value = 4
$ws.Range("M8").Value = "VCC, DATA, NC, GND"
$ws.Range("N8").Value = "Digital"
$ws.Range("O8").Value = "16 bits"
$ws.Range("P8").Value = "desde -40°C a 80°C"
$ws.Range("Q8").Value = 20
$ws.Range("R8").Value = "https://naylampmechatronics.com/sensores-temperatura-y-humedad/58-sensor-de-temperatura-y-humedad-relativa-dht22-am2302.html"

# Columnas D,E,F,G,H,I,J,K,L,N,O sin ajuste de texto (igual que fila 4)
$ws.Range("D8:L8").WrapText = $false
$ws.Range("N8:O8").WrapText = $false

# B8, C8, M8, P8, R8 mantienen ajuste de texto
$ws.Range("B8").WrapText = $true
$ws.Range("C8").WrapText = $true
$ws.Range("M8").WrapText = $true
$ws.Range("P8").WrapText = $true
$ws.Range("R8").WrapText = $true

# Q8 usa el mismo formato de moneda (USD) que el resto de la columna Precio
$ws.Range("Q8").NumberFormat = $ws.Range("Q7").NumberFormat
$ws.Range("Q8").WrapText = $true

# ---------------------------------------------------------------------
# 5) Ajustes de presentación de la hoja (anchos de columna, zoom, selección)
# ---------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 4.83
$ws.Columns("B").ColumnWidth = 27.3

$ws.Application.ActiveWindow.Zoom = 77
$ws.Range("K5").Select()
